# Prototype .csv final file generator — rework of the "tests" and "scope" sheets.
$wb = $excel.ActiveWorkbook

$wsTests = $wb.Worksheets.Item("tests")
$wsScope = $wb.Worksheets.Item("scope")

# ---------------------------------------------------------------------------
# "tests" sheet (sheet2): add an ID column in A, shift the old "ID" column
# from C to A, and introduce two new numeric columns "b" / "c" in B / C.
# Column D (True/False) keeps its existing TRUE/FALSE values.
# ---------------------------------------------------------------------------
$wsTests.Cells.Item(1,1).Value = "ID"
$wsTests.Cells.Item(1,2).Value = "b"
$wsTests.Cells.Item(1,3).Value = "c"
$wsTests.Cells.Item(1,4).Value = "True/False"

$testResults = @(
    @("test_1",  "TRUE"),
    @("test_2",  "TRUE"),
    @("test_3",  "TRUE"),
    @("test_4",  "TRUE"),
    @("test_5",  "FALSE"),
    @("test_6",  "FALSE"),
    @("test_7",  "TRUE"),
    @("test_8",  "TRUE"),
    @("test_9",  "TRUE"),
    @("test_10", "TRUE"),
    @("test_11", "TRUE"),
    @("test_12", "TRUE"),
    @("test_13", "TRUE"),
    @("test_14", "FALSE"),
    @("test_15", "TRUE"),
    @("test_16", "TRUE"),
    @("test_17", "TRUE"),
    @("test_18", "TRUE"),
    @("test_19", "FALSE"),
    @("test_20", "TRUE")
)

$r = 2
foreach ($row in $testResults) {
    $wsTests.Cells.Item($r,1).Value = $row[0]
    $wsTests.Cells.Item($r,2).Value = 0
    $wsTests.Cells.Item($r,3).Value = 0
    $wsTests.Cells.Item($r,4).Value = $row[1]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# "scope" sheet (sheet3): rebuild with the "basic" subset of tests first,
# followed by the full "extended" set of tests.
# ---------------------------------------------------------------------------
$basicTests = @(
    "test_1","test_2","test_4","test_5","test_7","test_8","test_9","test_10",
    "test_12","test_13","test_14","test_16","test_17","test_18","test_19","test_20"
)
$extendedTests = @(
    "test_1","test_2","test_3","test_4","test_5","test_6","test_7","test_8",
    "test_9","test_10","test_11","test_12","test_13","test_14","test_15",
    "test_16","test_17","test_18","test_19","test_20"
)

$wsScope.Cells.Item(1,1).Value = "ID"
$wsScope.Cells.Item(1,4).Value = "T/N"

$r = 2
foreach ($name in $basicTests) {
    $wsScope.Cells.Item($r,1).Value = $name
    $wsScope.Cells.Item($r,4).Value = "basic"
    $r = $r + 1
}
foreach ($name in $extendedTests) {
    $wsScope.Cells.Item($r,1).Value = $name
    $wsScope.Cells.Item($r,4).Value = "extended"
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Selection / active-sheet bookkeeping to match the saved view state.
# ---------------------------------------------------------------------------
$wsTests.Activate()
$wsTests.Range("H19").Select()

$wsScope.Range("D4").Select()
